$p = $ppt.ActivePresentation

function Restore-Geometry($shape, $origLeft, $origTop, $origWidth, $origHeight) {
    # The underlying text box uses <a:spAutoFit/>, so PowerPoint re-lays the
    # shape out whenever its text changes. Put the geometry back the way it
    # was (the canonical edit only touches run text) -- but only touch the
    # property if it actually drifted, since every COM round-trip through
    # points re-quantizes the EMU value by +/-1 and we don't want to
    # introduce needless noise on shapes that didn't really move.
    if ($shape.Left -ne $origLeft) { $shape.Left = $origLeft }
    if ($shape.Top -ne $origTop) { $shape.Top = $origTop }
    if ($shape.Width -ne $origWidth) { $shape.Width = $origWidth }
    if ($shape.Height -ne $origHeight) { $shape.Height = $origHeight }
}

function Fix-RunText($shape, $oldSnippet, $newSnippet) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    if ($full.IndexOf($oldSnippet) -lt 0) { return $false }

    $origLeft = $shape.Left
    $origTop = $shape.Top
    $origWidth = $shape.Width
    $origHeight = $shape.Height

    $tr.Replace($oldSnippet, $newSnippet, 1, 0, 0) | Out-Null

    Restore-Geometry $shape $origLeft $origTop $origWidth $origHeight
    return $true
}

function Merge-RunText($shape, $expectedText) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    if ($full -ne $expectedText) { return $false }

    $origLeft = $shape.Left
    $origTop = $shape.Top
    $origWidth = $shape.Width
    $origHeight = $shape.Height

    # Re-assigning the identical text over the whole range collapses the
    # (identically formatted) adjacent runs that made up that text into one.
    $sub = $tr.Characters(1, $full.Length)
    $sub.Text = $full

    Restore-Geometry $shape $origLeft $origTop $origWidth $origHeight
    return $true
}

# ---------------------------------------------------------------------------
# 1) & 2) Slide 2 ("Linear Probing Example"): two textboxes each contain a
#    "For example, quadratic probing for NN" line that should read "linear"
#    instead of "quadratic" (copy/paste leftover from the Quadratic slide).
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shape = $slide2.Shapes.Item($i)
    if ($shape.HasTextFrame -eq $false) { continue }
    Fix-RunText $shape "quadratic " "linear " | Out-Null
}

# ---------------------------------------------------------------------------
# 3) & 4) Slide 6 ("Double Hashing Example"): two adjacent runs that had been
#    split apart ("% " / "5" and "h(x) = x % " / "7") get merged back into a
#    single run each.
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shape = $slide6.Shapes.Item($i)
    if ($shape.HasTextFrame -eq $false) { continue }
    Merge-RunText $shape "h'(x) = x % 5" | Out-Null
    Merge-RunText $shape "h(x) = x % 7" | Out-Null
}
